# edit.ps1 - Reproduce the "Add files via upload" commit:
#   * remove the "thyroid" dataset column, keep "ecoli1" (which moves left
#     into the vacated column), and append four new dataset columns:
#     pima, yeast6, page_blocks, glass, each with fresh parameter values,
#     for every parameter sheet (SMOTE, DBSMOTE-eps, DBSMOTE-min,
#     WSSMOTE-nn, WSSMOTE-n_add, GSMOTE-nn).
#   * the previously-active sheet (WSSMOTE-n_add) is no longer the one
#     left selected; GSMOTE-nn (last sheet) ends up the active tab.

$wb = $excel.ActiveWorkbook

# xlPasteFormats
$xlPasteFormats = -4122

# New header labels for columns H:K (column G keeps "ecoli1", which is
# simply the old column H moved one slot to the left). The shared-string
# table records brand new strings in first-seen order, so to match the
# canonical file (where the new unique strings end up ordered
# yeast6, pima, page_blocks, glass) we must write column I before column H.
$newHeaders = @{ 'H' = 'pima'; 'I' = 'yeast6'; 'J' = 'page_blocks'; 'K' = 'glass' }
$headerWriteOrder = @('I', 'H', 'J', 'K')

# New parameter values (and, implicitly via copy-from-reference, styles)
# for columns H:K, rows 2-5, per sheet name.
$data = @{
    'SMOTE' = @{
        2 = @{ H = 2;   I = 10;  J = 5;   K = 5 }
        3 = @{ H = 5;   I = 2;   J = 2;   K = 5 }
        4 = @{ H = 2;   I = 10;  J = 2;   K = 5 }
        5 = @{ H = 2;   I = 20;  J = 2;   K = 5 }
    }
    'DBSMOTE-eps' = @{
        2 = @{ H = 0.8; I = 0.2; J = 0.2; K = 0.2 }
        3 = @{ H = 1.5; I = 0.2; J = 0.5; K = 1.2 }
        4 = @{ H = 0.5; I = 0.8; J = 0.2; K = 1 }
        5 = @{ H = 1.5; I = 1.2; J = 1;   K = 1.2 }
    }
    'DBSMOTE-min' = @{
        2 = @{ H = 3;   I = 3;   J = 2;   K = 2 }
        3 = @{ H = 2;   I = 5;   J = 5;   K = 3 }
        4 = @{ H = 2;   I = 5;   J = 2;   K = 5 }
        5 = @{ H = 2;   I = 5;   J = 3;   K = 2 }
    }
    'WSSMOTE-nn' = @{
        2 = @{ H = 10;  I = 5;   J = 2;   K = 5 }
        3 = @{ H = 5;   I = 5;   J = 10;  K = 5 }
        4 = @{ H = 20;  I = 2;   J = 5;   K = 5 }
        5 = @{ H = 10;  I = 20;  J = 2;   K = 2 }
    }
    'WSSMOTE-n_add' = @{
        2 = @{ H = 10;  I = 5;   J = 5;   K = 5 }
        3 = @{ H = 2;   I = 2;   J = 10;  K = 5 }
        4 = @{ H = 2;   I = 3;   J = 2;   K = 3 }
        5 = @{ H = 3;   I = 2;   J = 5;   K = 10 }
    }
    'GSMOTE-nn' = @{
        2 = @{ H = 20;  I = 10;  J = 5;   K = 5 }
        3 = @{ H = 2;   I = 2;   J = 10;  K = 2 }
        4 = @{ H = 20;  I = 20;  J = 10;  K = 5 }
        5 = @{ H = 10;  I = 20;  J = 2;   K = 5 }
    }
}

# Rows 2-5 in sheet "DBSMOTE-eps" (col H) and "DBSMOTE-min" (row 5, col G)
# need the borders/shading that match style index 6 in the original file;
# every other new data cell matches style index 5. Sheet "DBSMOTE-min" row5
# col G is produced automatically because it is a straight copy of the old
# column H (which already carried that particular style).
$style6Rows = @{
    'DBSMOTE-eps' = @(2, 5)
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $data.ContainsKey($name)) { continue }

    # --- Header row -----------------------------------------------------
    # Column G becomes "ecoli1": move the old column H header (value +
    # style) one slot to the left.
    $ws.Range("H1").Copy() | Out-Null
    $ws.Range("G1").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("G1").Value = $ws.Range("H1").Value2

    # Columns H:K get the four new dataset names, using the same header
    # style (already on H1).
    foreach ($col in $headerWriteOrder) {
        $ws.Range("H1").Copy() | Out-Null
        $ws.Range("$col`1").PasteSpecial($xlPasteFormats) | Out-Null
        $ws.Range("$col`1").Value = $newHeaders[$col]
    }

    # --- Data rows 2-5 ----------------------------------------------------
    foreach ($row in 2..5) {
        # Column G: move the old column H value+style one slot left
        # (this is literally the "ecoli1" column, unchanged).
        $ws.Range("H$row").Copy() | Out-Null
        $oldHValue = $ws.Range("H$row").Value2
        $ws.Range("G$row").PasteSpecial($xlPasteFormats) | Out-Null
        $ws.Range("G$row").Value = $oldHValue

        # Columns H:K: new parameter values for pima/yeast6/page_blocks/glass.
        $useStyle6 = $false
        if ($style6Rows.ContainsKey($name) -and $style6Rows[$name] -contains $row) {
            $useStyle6 = $true
        }

        foreach ($col in @('H', 'I', 'J', 'K')) {
            if ($useStyle6) {
                $wsRef = $wb.Worksheets.Item('DBSMOTE-min')
                $wsRef.Range("H5").Copy() | Out-Null
            } else {
                $wsRef = $wb.Worksheets.Item('SMOTE')
                $wsRef.Range("H2").Copy() | Out-Null
            }
            $ws.Range("$col$row").PasteSpecial($xlPasteFormats) | Out-Null
            $ws.Range("$col$row").Value = $data[$name][$row][$col]
        }
    }
}

# The workbook was left with the last sheet (GSMOTE-nn) active/selected.
$wb.Worksheets.Item('GSMOTE-nn').Activate()
$wb.Worksheets.Item('GSMOTE-nn').Range("A1").Select() | Out-Null
